$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (ras-a75-fw_005 -> ras-a75-fw_00F)
$ws.Name = "ras-a75-fw_00F"

# Clock signal renamed: i_clk -> o_clk (cell M22)
$ws.Range("M22").Value = "o_clk"

# Update the active selection to reflect the last edited cell
$ws.Range("M22").Select()
